$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-07-19 Saturday" "2025-07-20 Sunday"

Replace-Text "338÷2=169, 0" "504÷3=168, 0"
Replace-Text "549÷3=183, 0" "551÷8=68, 7"
Replace-Text "768÷3=256, 0" "285÷2=142, 1"
Replace-Text "565÷4=141, 1" "756÷9=84, 0"
Replace-Text "211÷9=23, 4" "794÷3=264, 2"
Replace-Text "309÷4=77, 1" "826÷8=103, 2"
Replace-Text "850÷7=121, 3" "254÷2=127, 0"
Replace-Text "801÷9=89, 0" "398÷8=49, 6"
Replace-Text "352÷8=44, 0" "222÷6=37, 0"
Replace-Text "968÷4=242, 0" "620÷2=310, 0"
Replace-Text "492÷3=164, 0" "499÷7=71, 2"
Replace-Text "412÷6=68, 4" "892÷6=148, 4"
Replace-Text "230÷3=76, 2" "277÷8=34, 5"
Replace-Text "913÷7=130, 3" "944÷2=472, 0"
Replace-Text "796÷4=199, 0" "383÷2=191, 1"
Replace-Text "811÷4=202, 3" "519÷9=57, 6"
Replace-Text "320÷6=53, 2" "227÷9=25, 2"
Replace-Text "493÷6=82, 1" "379÷8=47, 3"
Replace-Text "514÷6=85, 4" "468÷3=156, 0"
Replace-Text "658÷8=82, 2" "120÷4=30, 0"
Replace-Text "691÷6=115, 1" "530÷5=106, 0"
Replace-Text "675÷6=112, 3" "749÷4=187, 1"
Replace-Text "234÷6=39, 0" "612÷2=306, 0"
Replace-Text "860÷3=286, 2" "738÷3=246, 0"
Replace-Text "194÷5=38, 4" "979÷3=326, 1"
